$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three new shared-string backed rows (32, 33, 34) ------------
# Set values & formulas first so the dependency graph (e.g. F2's
# SUM(C:C)) picks them up correctly; copy formatting down afterwards.

$ws.Range("A32").Value = 42953
$ws.Range("C32").Value = 6
$ws.Range("D32").Value = "Adding test cases and persistance to QTableBrain."

$ws.Range("A33").Value = 42954
$ws.Range("C33").Value = 2.5
$ws.Range("D33").Value = "debugging persistance."

$ws.Range("A34").Value = 42955
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = "Refactored the choosing of random locations so always guaranteed to be in an empty cell."

# Fill the B column formula down across the new rows in one shot so it
# is recorded as a (valid) shared formula group, mirroring how dragging
# the fill handle down the "=A#" column works.
$ws.Range("B32:B34").Formula = "=A32"

# Copy the formatting (number formats / styles) from the prior row (31)
# down onto the three new rows without disturbing the values/formulas
# just written.
$ws.Range("A31:D31").Copy()
$ws.Range("A32:D34").PasteSpecial(-4122)

# --- Update the view state (selection / scroll position) -----------------
$excel.ActiveWindow.ScrollRow = 32
$ws.Range("D35").Select() | Out-Null
